# "Added more slides to the How section"
#
# Summary of edits:
#  - Slide 6 ("How - diversity of plaques in London"): rewrite body text,
#    drop the leading paragraph, reflow/resize the text box, drop the old
#    screenshot picture and instead pull in the chart graphic (Rectangle +
#    Picture group) that used to live on the next slide.
#  - The old slide 7 (chart recap slide) is removed - its chart group has
#    been relocated onto slide 6 and its caption text box is no longer
#    needed.
#  - The old slide 8 ("How - individual(s) with links to the slave trade")
#    becomes the new slide 7, with its title pluralised and its body
#    filled out with the full bullet list describing the scoring system.
#  - A near duplicate of that slide is added as slide 8 (a work-in-progress
#    slide the author had started).
#  - A brand new "Conclusion" slide is appended as slide 9.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 6 - rewrite the body placeholder text & reflow the shape
# ---------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$body6 = $s6.Shapes.Item(2)

$body6.TextFrame.TextRange.Text = "To retrieve the ethnicity of whose on the plaques, I used a qwikidata api to fetch data from the wikidata database for individual on the plaques.`rNot all entries had ethnicity fields as a property so the chart could be unreliable if there are correlations to the entry being there e.g. if the entry is more likely to be there if the individual is say White British.`rThe chart shows the overwhelming majority of those commemorated are British/English."

$body6.Left = 51.63842519685039
$body6.Top = 152.4252755905512
$body6.Width = 300.2986614173228
$body6.Height = 338.0
$body6.TextFrame.AutoSize = 2

# Drop the old screenshot picture that used to sit on slide 6
$s6.Shapes.Item(3).Delete()

# ---------------------------------------------------------------------
# 2. Pull the chart group over from the old slide 7 onto slide 6, then
#    remove the old slide 7.
# ---------------------------------------------------------------------
$s7old = $p.Slides.Item(7)
$chartGroup = $s7old.Shapes.Item(2)
$chartGroup.Copy()
$pasted = $s6.Shapes.Paste()
$grp6 = $pasted.Item(1)
$grp6.Left = 359.61535433070867
$grp6.Top = 152.4251968503937
$grp6.Width = 548.0
$grp6.Height = 338.0

$s7old.Delete()

# ---------------------------------------------------------------------
# 3. The old slide 8 is now at index 7. Update its title and body text.
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$title7 = $s7.Shapes.Item(1)
$title7.TextFrame.TextRange.Text = "How – individuals with links to the slave trade"

$body7 = $s7.Shapes.Item(2)
$body7.TextFrame.TextRange.Text = "To identify bad actors I made a scoring system taking into account the following:`rThe presence of certain words on their Wikipedia pages`rWhether they appear on the slave owner list on Wikipedia or Wikidata`rWhether they appear on the abolitionist list on Wikipedia`rEthnicity`rBirth period`rThe above criteria gave both positive and negative scores (higher the score the worse).`rWords such as `u2018poet`u2019/`u2019author`u2019/`u2019composer`u2019 appearing on the Wikipedia pages reduced the score, assuming the correlation that most poets/writers/`u2026 were against slavery.`r"

$body7.TextFrame.TextRange.Paragraphs(2,1).IndentLevel = 2
$body7.TextFrame.TextRange.Paragraphs(3,1).IndentLevel = 2
$body7.TextFrame.TextRange.Paragraphs(4,1).IndentLevel = 2
$body7.TextFrame.TextRange.Paragraphs(5,1).IndentLevel = 2
$body7.TextFrame.TextRange.Paragraphs(6,1).IndentLevel = 2
$body7.TextFrame.TextRange.Paragraphs(9,1).IndentLevel = 2

# ---------------------------------------------------------------------
# 4. Duplicate slide 7 to create slide 8 (a work-in-progress duplicate)
# ---------------------------------------------------------------------
$dup8 = $s7.Duplicate()
$s8 = $p.Slides.Item(8)
$body8 = $s8.Shapes.Item(2)
$body8.TextFrame.TextRange.Text = "The "

# ---------------------------------------------------------------------
# 5. Duplicate slide 8 to create slide 9, then turn it into "Conclusion"
# ---------------------------------------------------------------------
$dup9 = $s8.Duplicate()
$s9 = $p.Slides.Item(9)
$title9 = $s9.Shapes.Item(1)
$title9.TextFrame.TextRange.Text = "Conclusion"
$body9 = $s9.Shapes.Item(2)
$body9.TextFrame.TextRange.Text = ""
